$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.162809
$ws.Range("H2").Value = 6.488427000000001
$ws.Range("I2").Value = 0.06755089002018773
$ws.Range("J2").Value = 0.06755089002018773
$ws.Range("M2").Value = 154.942487
$ws.Range("N2").Value = 464.827461
$ws.Range("O2").Value = 0.982851703624775
$ws.Range("P2").Value = 0.9828517036247751
$ws.Range("Q2").Value = 335.111005365983
$ws.Range("R2").Value = 3015.999048293847
$ws.Range("S2").Value = 0.06639250733771132
$ws.Range("T2").Value = 0.06639250733771133
$ws.Range("G3").Value = 2.162809
$ws.Range("H3").Value = 6.488427000000001
$ws.Range("I3").Value = 0.06755089002018773
$ws.Range("J3").Value = 0.06755089002018773
$ws.Range("O3").Value = 0.003358739549735124
$ws.Range("P3").Value = 0.003358739549735124
$ws.Range("Q3").Value = 1.145188621155333
$ws.Range("R3").Value = 10.306697590398
$ws.Range("S3").Value = 0.0002268858459306122
$ws.Range("T3").Value = 0.0002268858459306122
$ws.Range("G4").Value = 2.162809
$ws.Range("H4").Value = 6.488427000000001
$ws.Range("I4").Value = 0.06755089002018773
$ws.Range("J4").Value = 0.06755089002018773
$ws.Range("M4").Value = 1.771368666666667
$ws.Range("N4").Value = 5.314106000000001
$ws.Range("O4").Value = 0.01123638032078883
$ws.Range("P4").Value = 0.01123638032078884
$ws.Range("Q4").Value = 3.831132094584667
$ws.Range("R4").Value = 34.48018885126201
$ws.Range("S4").Value = 0.0007590274912746082
$ws.Range("T4").Value = 0.0007590274912746084
$ws.Range("G5").Value = 2.162809
$ws.Range("H5").Value = 6.488427000000001
$ws.Range("I5").Value = 0.06755089002018773
$ws.Range("J5").Value = 0.06755089002018773
$ws.Range("M5").Value = 0.4024976666666666
$ws.Range("N5").Value = 1.207493
$ws.Range("O5").Value = 0.002553176504700935
$ws.Range("P5").Value = 0.002553176504700936
$ws.Range("Q5").Value = 0.8705255759456667
$ws.Range("R5").Value = 7.834730183511001
$ws.Range("S5").Value = 0.0001724693452711802
$ws.Range("T5").Value = 0.0001724693452711802
$ws.Range("I6").Value = 0.5628021396814664
$ws.Range("J6").Value = 0.5628021396814664
$ws.Range("M6").Value = 154.942487
$ws.Range("N6").Value = 464.827461
$ws.Range("O6").Value = 0.982851703624775
$ws.Range("P6").Value = 0.9828517036247751
$ws.Range("Q6").Value = 2791.986764266448
$ws.Range("R6").Value = 25127.88087839803
$ws.Range("S6").Value = 0.5531510417895978
$ws.Range("T6").Value = 0.553151041789598
$ws.Range("I7").Value = 0.5628021396814664
$ws.Range("J7").Value = 0.5628021396814664
$ws.Range("O7").Value = 0.003358739549735124
$ws.Range("P7").Value = 0.003358739549735124
$ws.Range("S7").Value = 0.001890305805223693
$ws.Range("T7").Value = 0.001890305805223693
$ws.Range("I8").Value = 0.5628021396814664
$ws.Range("J8").Value = 0.5628021396814664
$ws.Range("M8").Value = 1.771368666666667
$ws.Range("N8").Value = 5.314106000000001
$ws.Range("O8").Value = 0.01123638032078883
$ws.Range("P8").Value = 0.01123638032078884
$ws.Range("Q8").Value = 31.91918477447467
$ws.Range("R8").Value = 287.2726629702721
$ws.Range("S8").Value = 0.006323858886814678
$ws.Range("T8").Value = 0.006323858886814679
$ws.Range("I9").Value = 0.5628021396814664
$ws.Range("J9").Value = 0.5628021396814664
$ws.Range("M9").Value = 0.4024976666666666
$ws.Range("N9").Value = 1.207493
$ws.Range("O9").Value = 0.002553176504700935
$ws.Range("P9").Value = 0.002553176504700936
$ws.Range("Q9").Value = 7.252808314490666
$ws.Range("R9").Value = 65.275274830416
$ws.Range("S9").Value = 0.001436933199830134
$ws.Range("T9").Value = 0.001436933199830134
$ws.Range("G10").Value = 4.650307000000001
$ws.Range("H10").Value = 13.950921
$ws.Range("I10").Value = 0.1452427730405732
$ws.Range("J10").Value = 0.1452427730405732
$ws.Range("M10").Value = 154.942487
$ws.Range("N10").Value = 464.827461
$ws.Range("O10").Value = 0.982851703624775
$ws.Range("P10").Value = 0.9828517036247751
$ws.Range("Q10").Value = 720.5301318935091
$ws.Range("R10").Value = 6484.771187041581
$ws.Range("S10").Value = 0.142752106922114
$ws.Range("T10").Value = 0.142752106922114
$ws.Range("G11").Value = 4.650307000000001
$ws.Range("H11").Value = 13.950921
$ws.Range("I11").Value = 0.1452427730405732
$ws.Range("J11").Value = 0.1452427730405732
$ws.Range("O11").Value = 0.003358739549735124
$ws.Range("P11").Value = 0.003358739549735124
$ws.Range("Q11").Value = 2.462297253839334
$ws.Range("R11").Value = 22.160675284554
$ws.Range("S11").Value = 0.0004878326461245758
$ws.Range("T11").Value = 0.0004878326461245758
$ws.Range("G12").Value = 4.650307000000001
$ws.Range("H12").Value = 13.950921
$ws.Range("I12").Value = 0.1452427730405732
$ws.Range("J12").Value = 0.1452427730405732
$ws.Range("M12").Value = 1.771368666666667
$ws.Range("N12").Value = 5.314106000000001
$ws.Range("O12").Value = 0.01123638032078883
$ws.Range("P12").Value = 0.01123638032078884
$ws.Range("Q12").Value = 8.237408110180668
$ws.Range("R12").Value = 74.13667299162601
$ws.Range("S12").Value = 0.001632003036729896
$ws.Range("T12").Value = 0.001632003036729896
$ws.Range("G13").Value = 4.650307000000001
$ws.Range("H13").Value = 13.950921
$ws.Range("I13").Value = 0.1452427730405732
$ws.Range("J13").Value = 0.1452427730405732
$ws.Range("M13").Value = 0.4024976666666666
$ws.Range("N13").Value = 1.207493
$ws.Range("O13").Value = 0.002553176504700935
$ws.Range("P13").Value = 0.002553176504700936
$ws.Range("Q13").Value = 1.871737716783667
$ws.Range("R13").Value = 16.845639451053
$ws.Range("S13").Value = 0.000370830435604802
$ws.Range("T13").Value = 0.0003708304356048021
$ws.Range("G14").Value = 7.184856000000001
$ws.Range("H14").Value = 21.554568
$ws.Range("I14").Value = 0.2244041972577726
$ws.Range("J14").Value = 0.2244041972577726
$ws.Range("M14").Value = 154.942487
$ws.Range("N14").Value = 464.827461
$ws.Range("O14").Value = 0.982851703624775
$ws.Range("P14").Value = 0.9828517036247751
$ws.Range("Q14").Value = 1113.239457376872
$ws.Range("R14").Value = 10019.15511639185
$ws.Range("S14").Value = 0.2205560475753519
$ws.Range("T14").Value = 0.2205560475753519
$ws.Range("G15").Value = 7.184856000000001
$ws.Range("H15").Value = 21.554568
$ws.Range("I15").Value = 0.2244041972577726
$ws.Range("J15").Value = 0.2244041972577726
$ws.Range("O15").Value = 0.003358739549735124
$ws.Range("P15").Value = 0.003358739549735124
$ws.Range("Q15").Value = 3.804318983248
$ws.Range("R15").Value = 34.238870849232
$ws.Range("S15").Value = 0.0007537152524562432
$ws.Range("T15").Value = 0.0007537152524562432
$ws.Range("G16").Value = 7.184856000000001
$ws.Range("H16").Value = 21.554568
$ws.Range("I16").Value = 0.2244041972577726
$ws.Range("J16").Value = 0.2244041972577726
$ws.Range("M16").Value = 1.771368666666667
$ws.Range("N16").Value = 5.314106000000001
$ws.Range("O16").Value = 0.01123638032078883
$ws.Range("P16").Value = 0.01123638032078884
$ws.Range("Q16").Value = 12.727028792912
$ws.Range("R16").Value = 114.543259136208
$ws.Range("S16").Value = 0.002521490905969652
$ws.Range("T16").Value = 0.002521490905969652
$ws.Range("G17").Value = 7.184856000000001
$ws.Range("H17").Value = 21.554568
$ws.Range("I17").Value = 0.2244041972577726
$ws.Range("J17").Value = 0.2244041972577726
$ws.Range("M17").Value = 0.4024976666666666
$ws.Range("N17").Value = 1.207493
$ws.Range("O17").Value = 0.002553176504700935
$ws.Range("P17").Value = 0.002553176504700936
$ws.Range("Q17").Value = 2.891887775336
$ws.Range("R17").Value = 26.026989978024
$ws.Range("S17").Value = 0.0005729435239948191
$ws.Range("T17").Value = 0.0005729435239948192
